# Update the "YOUR TEAM" contact slide (last slide) with new contact info:
# Lee Harrison's entry is replaced with Carrie Rotman's, Bernie's phone
# extension / email casing is corrected, and the company phone/domain are
# refreshed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

$bullet = [char]0x2022

# Lee Harrison -> Carrie Rotman
$s.Shapes.Item("Text 6").TextFrame.TextRange.Text = "Carrie Rotman"

# M: 289-388-8399 -> O: 1-877-660-3660 ext. 235
$s.Shapes.Item("Text 7").TextFrame.TextRange.Text = "O: 1-877-660-3660 ext. 235"

# lee@allianceglobaladvisors.com -> carrie@alliancefinancing.com
$s.Shapes.Item("Text 8").TextFrame.TextRange.Text = "carrie@alliancefinancing.com"

# O: 905-660-3660 ext 225   .   M: 416-569-2899 -> O: 1-877-660-3660 ext. 222   .   M: 416-569-2899
$s.Shapes.Item("Text 12").TextFrame.TextRange.Text = "O: 1-877-660-3660 ext. 222   $bullet   M: 416-569-2899"

# Bernie@alliancefinancing.com -> bernie@alliancefinancing.com
$s.Shapes.Item("Text 13").TextFrame.TextRange.Text = "bernie@alliancefinancing.com"

# alliancefinancing.com -> alliancefinancing.ai
$s.Shapes.Item("Text 16").TextFrame.TextRange.Text = "alliancefinancing.ai"

# 905-660-3660 -> 1-877-660-3660
$s.Shapes.Item("Text 17").TextFrame.TextRange.Text = "1-877-660-3660"
